{"js": "const replacements = [\n  [\"2025-09-15 Monday\", \"2025-09-16 Tuesday\"],\n  [\"209\u00f79=\", \"534\u00f78=\"],\n  [\"721\u00f73=\", \"457\u00f75=\"],\n  [\"167\u00f73=\", \"367\u00f75=\"],\n  [\"255\u00f77=\", \"934\u00f72=\"],\n  [\"470\u00f72=\", \"316\u00f76=\"],\n  [\"618\u00f79=\", \"203\u00f76=\"],\n  [\"591\u00f79=\", \"222\u00f76=\"],\n  [\"900\u00f74=\", \"664\u00f74=\"],\n  [\"110\u00f78=\", \"574\u00f77=\"],\n  [\"103\u00f72=\", \"980\u00f76=\"],\n  [\"939\u00f72=\", \"626\u00f74=\"],\n  [\"677\u00f76=\", \"410\u00f76=\"],\n  [\"534\u00f73=\", \"341\u00f77=\"],\n  [\"547\u00f75=\", \"196\u00f75=\"],\n  [\"630\u00f79=\", \"453\u00f73=\"],\n  [\"517\u00f72=\", \"616\u00f72=\"],\n  [\"791\u00f74=\", \"519\u00f78=\"],\n  [\"327\u00f76=\", \"743\u00f76=\"],\n  [\"191\u00f72=\", \"645\u00f76=\"],\n  [\"387\u00f72=\", \"279\u00f76=\"],\n  [\"938\u00f75=\", \"449\u00f79=\"],\n  [\"728\u00f78=\", \"945\u00f76=\"],\n  [\"659\u00f75=\", \"840\u00f79=\"],\n  [\"339\u00f75=\", \"727\u00f74=\"],\n  [\"107\u00f74=\", \"183\u00f77=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-09-15 Monday', '2025-09-16 Tuesday'),\n    @('209\u00f79=', '534\u00f78='),\n    @('721\u00f73=', '457\u00f75='),\n    @('167\u00f73=', '367\u00f75='),\n    @('255\u00f77=', '934\u00f72='),\n    @('470\u00f72=', '316\u00f76='),\n    @('618\u00f79=', '203\u00f76='),\n    @('591\u00f79=', '222\u00f76='),\n    @('900\u00f74=', '664\u00f74='),\n    @('110\u00f78=', '574\u00f77='),\n    @('103\u00f72=', '980\u00f76='),\n    @('939\u00f72=', '626\u00f74='),\n    @('677\u00f76=', '410\u00f76='),\n    @('534\u00f73=', '341\u00f77='),\n    @('547\u00f75=', '196\u00f75='),\n    @('630\u00f79=', '453\u00f73='),\n    @('517\u00f72=', '616\u00f72='),\n    @('791\u00f74=', '519\u00f78='),\n    @('327\u00f76=', '743\u00f76='),\n    @('191\u00f72=', '645\u00f76='),\n    @('387\u00f72=', '279\u00f76='),\n    @('938\u00f75=', '449\u00f79='),\n    @('728\u00f78=', '945\u00f76='),\n    @('659\u00f75=', '840\u00f79='),\n    @('339\u00f75=', '727\u00f74='),\n    @('107\u00f74=', '183\u00f77=')\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.Execute(\n        $find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2\n    )\n}\n"}
